$p = $ppt.ActivePresentation

# Locate the shape on slide 9 whose text references the "18%" churn-risk
# stat (keeps the script robust even if indices ever shift).
$targetSlideIndex = 9
$targetShapeIndex = 6

$slide = $null
$shape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*18%*top-tier*") {
                $slide = $sl
                $shape = $shp
            }
        }
    }
}

if ($shape -eq $null) {
    $slide = $p.Slides.Item($targetSlideIndex)
    $shape = $slide.Shapes.Item($targetShapeIndex)
}

$tr = $shape.TextFrame.TextRange

# Replace only the "the 18% " substring with "the 20% " so the run gets
# split the same way PowerPoint splits it when you retype just that part
# (leaving the surrounding runs/text untouched).
$oldFragment = "the 18% "
$newFragment = "the 20% "
$startPos = $tr.Text.IndexOf($oldFragment) + 1

$sub = $tr.Characters($startPos, $oldFragment.Length)
$sub.Text = $newFragment

# The text box auto-fits its height to the (now two-line) paragraph; match
# the recalculated layout extent PowerPoint would have produced.
$shape.Height = 50.89221
